$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commodity_list")

# Populate the new commodity rows. Values are written in the same order
# the original author entered them, so shared-string de-dup indices come
# out in the same sequence as the target workbook.

# Row 3: Cocoa futures
$ws.Range("A3").Value = "CC=F"
$ws.Range("B3").Value = "Cocoa"

# Row 7: West Texas Intermediate crude futures
$ws.Range("A7").Value = "CL=F"
$ws.Range("B7").Value = "West Texas Intermediate"

# Row 2: Brent Crude futures (name entered before ticker)
$ws.Range("B2").Value = "Brent Crude"
$ws.Range("A2").Value = "BZ=F"

# Row 6: Natural Gas futures (name entered before ticker)
$ws.Range("B6").Value = "Natural Gas"
$ws.Range("A6").Value = "NG=F"

# Row 4: Gold futures (name entered before ticker)
$ws.Range("B4").Value = "Gold"
$ws.Range("A4").Value = "GC=F"

# Row 5: Silver futures (name entered before ticker)
$ws.Range("B5").Value = "Silver"
$ws.Range("A5").Value = "SI=F"

# Make Commodity_list the active sheet/tab and set its selection, moving the
# "tabSelected" flag off Treasury_list (which held it before this edit).
$ws.Activate()
$ws.Range("A6").Select()
